$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing the existing rows 3-4 down to 4-5.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44463
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 100112012
$ws.Range("G3").Value = "Espinaca"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = "`$/cuna 10 kilos"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 1200
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = "Hortaliza"

# Ensure the date cell keeps the date-number format used by the other date cells.
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
